# Regenerate the "K" column (column G, header "K") values in the save_data
# sheet. The values were recalculated from a K (strikeouts) source instead
# of the previous "Strike#" derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-32 (row index => value), per the regenerated data.
$kValues = @{
    2  = 3
    3  = 1
    4  = 2
    5  = 0
    6  = 2
    7  = 3
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 2
    13 = 1
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 2
    28 = 0
    29 = 5
    30 = 6
    31 = 3
    32 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
